# Drop the oldest day (2025-11-04, row 2) from the GSC export "Chart" sheet.
# Deleting the entire row shifts the remaining dates/values up by one and
# re-packs the shared-string table; the "Table" sheet's header row simply
# points at the same (now renumbered) shared strings automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows.Item(2).Delete()
